$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '26.715.08'
$ws.Range('E2').Value = '  +0.03%  '
Set-TextCell 'D3' '1.636.40'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('E4').Value = '  +0.26%  '
Set-TextCell 'D5' '217.16'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('E8').Value = '  -0.72%  '
Set-TextCell 'D9' '0.0621'
$ws.Range('E9').Value = '  -0.73%  '
Set-TextCell 'D10' '19.03'
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('E12').Value = '  -0.62%  '
Set-TextCell 'D13' '1.636.40'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E15').Value = '  -1.37%  '
Set-TextCell 'D16' '64.32'
$ws.Range('E16').Value = '  -1.32%  '
Set-TextCell 'D17' '26.701.91'
$ws.Range('E18').Value = '  -2.35%  '
Set-TextCell 'D20' '210.21'
$ws.Range('E20').Value = '  -3.52%  '
$ws.Range('E21').Value = '  -0.84%  '
Set-TextCell 'D22' '6.17'
$ws.Range('E22').Value = '  -1.50%  '
Set-TextCell 'D23' '2.32'
$ws.Range('E23').Value = '  +2.64%  '
$ws.Range('E24').Value = '  -2.99%  '
Set-TextCell 'D25' '145.56'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E27').Value = '  -2.44%  '
Set-TextCell 'D28' '7.05'
$ws.Range('E28').Value = '  -0.95%  '
Set-TextCell 'D29' '15.53'
$ws.Range('E29').Value = '  -1.17%  '
Set-TextCell 'D30' '0.0503'
$ws.Range('E30').Value = '  -2.43%  '
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('E33').Value = '  -1.43%  '
Set-TextCell 'D34' '1.274.23'
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('E35').Value = '  -1.41%  '
$ws.Range('E36').Value = '  +0.44%  '
$ws.Range('E37').Value = '  -1.99%  '
Set-TextCell 'D38' '0.529'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('E39').Value = '  -2.10%  '
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('E41').Value = '  -1.46%  '
$ws.Range('E42').Value = '  -2.20%  '
Set-TextCell 'D43' '1.774.00'
$ws.Range('E43').Value = '  -0.60%  '
$ws.Range('E44').Value = '  -3.49%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D45' '60.43'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D46' '91.03'
$ws.Range('E46').Value = '  -0.99%  '
$ws.Range('E47').Value = '  -2.43%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D48' '0.0519'
$ws.Range('E48').Value = '  +0.81%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D49' '7.53'
$ws.Range('E49').Value = '  -3.09%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D50' '0.0960'
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D51' '0.406'
$ws.Range('E51').Value = '  -0.24%  '
